$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 1101
$ws.Range("I33").Value = 323.33334
$ws.Range("K33").Value = 323.33334
$ws.Range("M33").Value = -94.33334000000002

$ws.Range("H40").Value = 928.1875
$ws.Range("I40").Value = 750.1111
$ws.Range("J40").Value = 1157.1428
$ws.Range("K40").Value = 750.1111
$ws.Range("L40").Value = 1157.1428
$ws.Range("M40").Value = -575.1111
$ws.Range("N40").Value = -1507.1428

$ws.Range("H49").Value = 642.5
$ws.Range("J49").Value = 770
$ws.Range("L49").Value = 2310
$ws.Range("N49").Value = -2582

$ws.Range("H64").Value = 3953.3333
$ws.Range("I64").Value = 2895
$ws.Range("J64").Value = 4800
$ws.Range("K64").Value = 2895
$ws.Range("L64").Value = 4800
$ws.Range("M64").Value = -2647
$ws.Range("N64").Value = -5296

$ws.Range("H67").Value = 3953.3333
$ws.Range("I67").Value = 2895
$ws.Range("J67").Value = 4800
$ws.Range("K67").Value = 2895
$ws.Range("L67").Value = 4800
$ws.Range("M67").Value = -2037
$ws.Range("N67").Value = -6516

$ws.Range("H76").Value = 3475091
$ws.Range("I76").Value = 3082.7144
$ws.Range("J76").Value = 6175542
$ws.Range("K76").Value = 3082.7144
$ws.Range("L76").Value = 6175542
$ws.Range("M76").Value = -2767.7144
$ws.Range("N76").Value = -6176172

$ws.Range("H79").Value = 3475091
$ws.Range("I79").Value = 3082.7144
$ws.Range("J79").Value = 6175542
$ws.Range("K79").Value = 3082.7144
$ws.Range("L79").Value = 6175542
$ws.Range("M79").Value = -1990.7144
$ws.Range("N79").Value = -6177726

$ws.Range("H80").Value = 12183507
$ws.Range("I80").Value = 500.25
$ws.Range("J80").Value = 20305512
$ws.Range("K80").Value = 1500.75
$ws.Range("L80").Value = 60916536
$ws.Range("M80").Value = -502.75
$ws.Range("N80").Value = -60918532

$ws.Range("H83").Value = 12183507
$ws.Range("I83").Value = 500.25
$ws.Range("J83").Value = 20305512
$ws.Range("K83").Value = 4502.25
$ws.Range("L83").Value = 182749608
$ws.Range("M83").Value = 489.75
$ws.Range("N83").Value = -182759592

$ws.Range("H113").Value = 111116370
$ws.Range("I113").Value = 1000000000
$ws.Range("J113").Value = 5913.375
$ws.Range("K113").Value = 1000000000
$ws.Range("L113").Value = 5913.375
$ws.Range("M113").Value = -999996746
$ws.Range("N113").Value = -12421.375

$ws.Range("H116").Value = 11765444
$ws.Range("I116").Value = 40324508
$ws.Range("J116").Value = 5829.5293
$ws.Range("K116").Value = 40324508
$ws.Range("L116").Value = 5829.5293
$ws.Range("M116").Value = -40321066
$ws.Range("N116").Value = -12713.5293

$ws.Range("H138").Value = 3046.1155
$ws.Range("I138").Value = 2245.2632
$ws.Range("J138").Value = 3507.2122
$ws.Range("K138").Value = 6735.7896
$ws.Range("L138").Value = 10521.6366
$ws.Range("M138").Value = -1595.7896
$ws.Range("N138").Value = -20801.6366

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2922.862
$ws.Range("I45").Value = 2424.318
$ws.Range("K45").Value = 2424.318
$ws.Range("M45").Value = -2047.318

$ws.Range("H61").Value = 12350039
$ws.Range("I61").Value = 15155680
$ws.Range("J61").Value = 5220
$ws.Range("K61").Value = 15155680
$ws.Range("L61").Value = 5220
$ws.Range("M61").Value = -15155468
$ws.Range("N61").Value = -5644

$ws.Range("H63").Value = 10417333

$ws.Range("H66").Value = 10417333

$ws.Range("H74").Value = 28572998
$ws.Range("I74").Value = 37037704
$ws.Range("J74").Value = 4612.5
$ws.Range("K74").Value = 37037704
$ws.Range("L74").Value = 4612.5
$ws.Range("M74").Value = -37036830
$ws.Range("N74").Value = -6360.5

$ws.Range("H77").Value = 28572998
$ws.Range("I77").Value = 37037704
$ws.Range("J77").Value = 4612.5
$ws.Range("K77").Value = 185188520
$ws.Range("L77").Value = 23062.5
$ws.Range("M77").Value = -185184152
$ws.Range("N77").Value = -31798.5

$ws.Range("H88").Value = 112640.664
$ws.Range("I88").Value = 1549.75
$ws.Range("J88").Value = 201513.4
$ws.Range("K88").Value = 1549.75
$ws.Range("L88").Value = 201513.4
$ws.Range("M88").Value = -1143.75
$ws.Range("N88").Value = -202325.4

$ws.Range("H91").Value = 112640.664
$ws.Range("I91").Value = 1549.75
$ws.Range("J91").Value = 201513.4
$ws.Range("K91").Value = 1549.75
$ws.Range("L91").Value = 201513.4
$ws.Range("M91").Value = -145.75
$ws.Range("N91").Value = -204321.4

$ws.Range("H128").Value = 39999.668
$ws.Range("J128").Value = 39999.668
$ws.Range("L128").Value = 39999.668
$ws.Range("N128").Value = -49959.668

$ws.Range("H132").Value = 10883150
$ws.Range("I132").Value = 15627323
$ws.Range("J132").Value = 39324.285
$ws.Range("K132").Value = 46881969
$ws.Range("L132").Value = 117972.855
$ws.Range("M132").Value = -46879439
$ws.Range("N132").Value = -123032.855

$ws.Range("H136").Value = 12350039
$ws.Range("I136").Value = 15155680
$ws.Range("J136").Value = 5220
$ws.Range("K136").Value = 45467040
$ws.Range("L136").Value = 15660
$ws.Range("M136").Value = -45464490
$ws.Range("N136").Value = -20760

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3500.5833
$ws.Range("I20").Value = 2874.75
$ws.Range("K20").Value = 2874.75
$ws.Range("M20").Value = -2627.75

$ws.Range("H22").Value = 972.5454999999999
$ws.Range("I22").Value = 966
$ws.Range("K22").Value = 966
$ws.Range("M22").Value = -793

$ws.Range("H105").Value = 6192642.5
$ws.Range("I105").Value = 15153063
$ws.Range("J105").Value = 2085783.2
$ws.Range("K105").Value = 15153063
$ws.Range("L105").Value = 2085783.2
$ws.Range("M105").Value = -15151316
$ws.Range("N105").Value = -2089277.2

$ws.Range("H134").Value = 7390.6816
$ws.Range("I134").Value = 8220.467000000001
$ws.Range("K134").Value = 24661.401
$ws.Range("M134").Value = -22126.401

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1519.8572
$ws.Range("I16").Value = 1606.5
$ws.Range("J16").Value = 1000
$ws.Range("K16").Value = 1606.5
$ws.Range("L16").Value = 1000
$ws.Range("M16").Value = -1319.5
$ws.Range("N16").Value = -1574

$ws.Range("H31").Value = 4210.3623
$ws.Range("I31").Value = 1809.5518
$ws.Range("J31").Value = 6611.1724
$ws.Range("K31").Value = 1809.5518
$ws.Range("L31").Value = 6611.1724
$ws.Range("M31").Value = -1514.5518
$ws.Range("N31").Value = -7201.1724

$ws.Range("H34").Value = 4210.3623
$ws.Range("I34").Value = 1809.5518
$ws.Range("J34").Value = 6611.1724
$ws.Range("K34").Value = 1809.5518
$ws.Range("L34").Value = 6611.1724
$ws.Range("M34").Value = -1607.5518
$ws.Range("N34").Value = -7015.1724

$ws.Range("H113").Value = 1519.8572
$ws.Range("I113").Value = 1606.5
$ws.Range("J113").Value = 1000
$ws.Range("K113").Value = 1606.5
$ws.Range("L113").Value = 1000
$ws.Range("M113").Value = 563.5
$ws.Range("N113").Value = -5340

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 50104000
$ws.Range("J37").Value = 50104000
$ws.Range("L37").Value = 150312000
$ws.Range("N37").Value = -150312224

$ws.Range("H95").Value = 3518

$ws.Range("H122").Value = 1378.9615
$ws.Range("I122").Value = 104
$ws.Range("K122").Value = 936
$ws.Range("M122").Value = 1514

$ws.Range("H131").Value = 680.24
$ws.Range("J131").Value = 714.2023
$ws.Range("L131").Value = 2142.6069
$ws.Range("N131").Value = -12222.6069

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 3680920.8
$ws.Range("I70").Value = 4771.2856
$ws.Range("K70").Value = 4771.2856
$ws.Range("M70").Value = -4501.2856

$ws.Range("H73").Value = 3680920.8
$ws.Range("I73").Value = 4771.2856
$ws.Range("K73").Value = 4771.2856
$ws.Range("M73").Value = -3835.2856

$ws.Range("H95").Value = 14260
$ws.Range("J95").Value = 14260
$ws.Range("L95").Value = 14260
$ws.Range("N95").Value = -19752

$ws.Range("H122").Value = 4105.069
$ws.Range("I122").Value = 3833.85
$ws.Range("J122").Value = 4707.778
$ws.Range("K122").Value = 11501.55
$ws.Range("L122").Value = 14123.334
$ws.Range("M122").Value = -9051.549999999999
$ws.Range("N122").Value = -19023.334

$ws.Range("H132").Value = 7505907.5
$ws.Range("I132").Value = 14118557
$ws.Range("J132").Value = 66676.5
$ws.Range("K132").Value = 42355671
$ws.Range("L132").Value = 200029.5
$ws.Range("M132").Value = -42353141
$ws.Range("N132").Value = -205089.5

$ws.Range("H140").Value = 56392
$ws.Range("J140").Value = 56392
$ws.Range("L140").Value = 56392
$ws.Range("N140").Value = -66752

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2437.5
$ws.Range("I68").Value = 2500
$ws.Range("J68").Value = 2375
$ws.Range("K68").Value = 2500
$ws.Range("L68").Value = 2375
$ws.Range("M68").Value = -1751
$ws.Range("N68").Value = -3873

$ws.Range("H71").Value = 2437.5
$ws.Range("I71").Value = 2500
$ws.Range("J71").Value = 2375
$ws.Range("K71").Value = 12500
$ws.Range("L71").Value = 11875
$ws.Range("M71").Value = -8756
$ws.Range("N71").Value = -19363

$ws.Range("H82").Value = 2733.8572
$ws.Range("I82").Value = 2297.6365
$ws.Range("K82").Value = 2297.6365
$ws.Range("M82").Value = -1936.6365

$ws.Range("H85").Value = 2733.8572
$ws.Range("I85").Value = 2297.6365
$ws.Range("K85").Value = 2297.6365
$ws.Range("M85").Value = -1049.6365

$ws.Range("H93").Value = 1777.7222
$ws.Range("I93").Value = 1799.9333
$ws.Range("K93").Value = 1799.9333
$ws.Range("M93").Value = -551.9332999999999

$ws.Range("H122").Value = 1157211.2
$ws.Range("I122").Value = 1403863.6
$ws.Range("K122").Value = 4211590.800000001
$ws.Range("M122").Value = -4209140.800000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H120").Value = 40000
$ws.Range("J120").Value = 40000
$ws.Range("L120").Value = 40000
$ws.Range("N120").Value = -49676

$ws.Range("H131").Value = 22083.334
$ws.Range("J131").Value = 22083.334
$ws.Range("L131").Value = 22083.334
$ws.Range("N131").Value = -32163.334
